$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
$ws.Range("H58").Value = 1589.4706
$ws.Range("J58").Value = 2554.111
$ws.Range("L58").Value = 7662.333
$ws.Range("N58").Value = -7962.333
$ws.Range("H74").Value = 153675
$ws.Range("I74").Value = 153675
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 153675
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -152739
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 153675
$ws.Range("I77").Value = 153675
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 768375
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -763695
$ws.Range("N77").ClearContents()
$ws.Range("H100").Value = 5383.385
$ws.Range("I100").Value = 4335
$ws.Range("J100").Value = 5697.9
$ws.Range("K100").Value = 4335
$ws.Range("L100").Value = 5697.9
$ws.Range("M100").Value = -3794
$ws.Range("N100").Value = -6779.9
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 5250.5
$ws.Range("I5").Value = 10001
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 10001
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -9889
$ws.Range("N5").Value = -724
$ws.Range("H63").Value = 4133.5
$ws.Range("I63").Value = 3759
$ws.Range("K63").Value = 3759
$ws.Range("M63").Value = -3073
$ws.Range("H66").Value = 4133.5
$ws.Range("I66").Value = 3759
$ws.Range("K66").Value = 18795
$ws.Range("M66").Value = -15363
$ws.Range("H88").Value = 3748
$ws.Range("I88").Value = 3996
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 3996
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -3590
$ws.Range("N88").Value = -4312
$ws.Range("H91").Value = 3748
$ws.Range("I91").Value = 3996
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 3996
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -2592
$ws.Range("N91").Value = -6308
$ws.Range("H101").Value = 43520.2
$ws.Range("J101").Value = 43520.2
$ws.Range("L101").Value = 43520.2
$ws.Range("N101").Value = -50010.2
$ws.Range("H102").Value = 6739.8
$ws.Range("I102").Value = 6739.8
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6739.8
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -5117.8
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 4962.2
$ws.Range("I122").Value = 5879.4546
$ws.Range("J122").Value = 2439.75
$ws.Range("K122").Value = 17638.3638
$ws.Range("L122").Value = 7319.25
$ws.Range("M122").Value = -15188.3638
$ws.Range("N122").Value = -12219.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 5250.5
$ws.Range("I4").Value = 10001
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 10001
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -9886
$ws.Range("N4").Value = -730
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H82").Value = 65025.7
$ws.Range("I82").Value = 65025.7
$ws.Range("K82").Value = 65025.7
$ws.Range("M82").Value = -64642.7
$ws.Range("H85").Value = 65025.7
$ws.Range("I85").Value = 65025.7
$ws.Range("K85").Value = 65025.7
$ws.Range("M85").Value = -63699.7
$ws.Range("H86").Value = 37101.414
$ws.Range("I86").Value = 2419.1765
$ws.Range("J86").Value = 86234.586
$ws.Range("K86").Value = 2419.1765
$ws.Range("L86").Value = 86234.586
$ws.Range("M86").Value = -1296.1765
$ws.Range("N86").Value = -88480.586
$ws.Range("H89").Value = 37101.414
$ws.Range("I89").Value = 2419.1765
$ws.Range("J89").Value = 86234.586
$ws.Range("K89").Value = 12095.8825
$ws.Range("L89").Value = 431172.93
$ws.Range("M89").Value = -6479.8825
$ws.Range("N89").Value = -442404.93
$ws.Range("H97").Value = 65248.445
$ws.Range("I97").Value = 67220.625
$ws.Range("K97").Value = 67220.625
$ws.Range("M97").Value = -66229.625
$ws.Range("H105").Value = 3375.5
$ws.Range("I105").Value = 3334.5
$ws.Range("J105").Value = 3498.5
$ws.Range("K105").Value = 3334.5
$ws.Range("L105").Value = 3498.5
$ws.Range("M105").Value = -1587.5
$ws.Range("N105").Value = -6992.5
$ws.Range("H132").Value = 90961.42999999999
$ws.Range("J132").Value = 90961.42999999999
$ws.Range("L132").Value = 90961.42999999999
$ws.Range("N132").Value = -101081.43

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 28987.5
$ws.Range("J109").Value = 28987.5
$ws.Range("L109").Value = 28987.5
$ws.Range("N109").Value = -31067.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 182923.7
$ws.Range("J68").Value = 1640
$ws.Range("L68").Value = 4920
$ws.Range("N68").Value = -6542
$ws.Range("H71").Value = 182923.7
$ws.Range("J71").Value = 1640
$ws.Range("L71").Value = 14760
$ws.Range("N71").Value = -22872
$ws.Range("H122").Value = 723.86664
$ws.Range("I122").Value = 527.61536
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 4748.53824
$ws.Range("L122").Value = 17995.5
$ws.Range("M122").Value = -2298.53824
$ws.Range("N122").Value = -22895.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3209.6667
$ws.Range("I122").Value = 3151.1667
$ws.Range("K122").Value = 9453.500100000001
$ws.Range("M122").Value = -7003.500100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 50000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H133").Value = 49800
$ws.Range("J133").Value = 49800
$ws.Range("L133").Value = 49800
$ws.Range("N133").Value = -54860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 33000
$ws.Range("J93").Value = 33000
$ws.Range("L93").Value = 33000
$ws.Range("N93").Value = -37992
